$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1716.6666
$ws.Range("J7").Value = 1575
$ws.Range("L7").Value = 1575
$ws.Range("N7").Value = -1799

$ws.Range("H14").Value = 1716.6666
$ws.Range("J14").Value = 1575
$ws.Range("L14").Value = 1575
$ws.Range("N14").Value = -1957

$ws.Range("H111").Value = 6576
$ws.Range("I111").Value = 6576
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 19728
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -16661
$ws.Range("N111").ClearContents()

$ws.Range("H118").Value = 2460.8215
$ws.Range("J118").Value = 3580.2856
$ws.Range("L118").Value = 10740.8568
$ws.Range("N118").Value = -14054.8568

$ws.Range("H127").Value = 399212.53
$ws.Range("I127").Value = 512.5
$ws.Range("J127").Value = 468551.66
$ws.Range("K127").Value = 1537.5
$ws.Range("L127").Value = 1405654.98
$ws.Range("M127").Value = 3422.5
$ws.Range("N127").Value = -1415574.98

$ws.Range("H129").Value = 994.7347
$ws.Range("J129").Value = 1009.68085
$ws.Range("L129").Value = 3029.04255
$ws.Range("N129").Value = -13029.04255

$ws.Range("H132").Value = 9625206
$ws.Range("I132").Value = 8738.412
$ws.Range("J132").Value = 27789644
$ws.Range("K132").Value = 26215.236
$ws.Range("L132").Value = 83368932
$ws.Range("M132").Value = -23685.236
$ws.Range("N132").Value = -83373992

$ws.Range("H137").Value = 6063257.5
$ws.Range("I137").Value = 1909.1364
$ws.Range("J137").Value = 18185954
$ws.Range("K137").Value = 5727.4092
$ws.Range("L137").Value = 54557862
$ws.Range("M137").Value = -3177.4092
$ws.Range("N137").Value = -54562962

$ws.Range("H138").Value = 6413287
$ws.Range("I138").Value = 1732.174
$ws.Range("J138").Value = 15629897
$ws.Range("K138").Value = 5196.522
$ws.Range("L138").Value = 46889691
$ws.Range("M138").Value = -56.52199999999993
$ws.Range("N138").Value = -46899971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6941.7896
$ws.Range("I32").Value = 6991.469
$ws.Range("K32").Value = 6991.469
$ws.Range("M32").Value = -6704.469

$ws.Range("H88").Value = 2608.1428
$ws.Range("I88").Value = 2333.3333
$ws.Range("K88").Value = 2333.3333
$ws.Range("M88").Value = -1927.3333

$ws.Range("H91").Value = 2608.1428
$ws.Range("I91").Value = 2333.3333
$ws.Range("K91").Value = 2333.3333
$ws.Range("M91").Value = -929.3332999999998

$ws.Range("H132").Value = 14709685
$ws.Range("I132").Value = 25003544
$ws.Range("J132").Value = 4173.143
$ws.Range("K132").Value = 75010632
$ws.Range("L132").Value = 12519.429
$ws.Range("M132").Value = -75008102
$ws.Range("N132").Value = -17579.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22729106
$ws.Range("I86").Value = 1899.5555
$ws.Range("J86").Value = 38463324
$ws.Range("K86").Value = 1899.5555
$ws.Range("L86").Value = 38463324
$ws.Range("M86").Value = -776.5554999999999
$ws.Range("N86").Value = -38465570

$ws.Range("H89").Value = 22729106
$ws.Range("I89").Value = 1899.5555
$ws.Range("J89").Value = 38463324
$ws.Range("K89").Value = 9497.7775
$ws.Range("L89").Value = 192316620
$ws.Range("M89").Value = -3881.7775
$ws.Range("N89").Value = -192327852

$ws.Range("H99").Value = 1607.826
$ws.Range("I99").Value = 1538.6666
$ws.Range("J99").Value = 1737.5
$ws.Range("K99").Value = 1538.6666
$ws.Range("L99").Value = 1737.5
$ws.Range("M99").Value = -40.66660000000002
$ws.Range("N99").Value = -4733.5

$ws.Range("H105").Value = 4294.4443
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 4456.25
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 4456.25
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -7950.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6176592.5
$ws.Range("I31").Value = 5403.0303
$ws.Range("J31").Value = 15874177
$ws.Range("K31").Value = 5403.0303
$ws.Range("L31").Value = 15874177
$ws.Range("M31").Value = -5108.0303
$ws.Range("N31").Value = -15874767

$ws.Range("H34").Value = 6176592.5
$ws.Range("I34").Value = 5403.0303
$ws.Range("J34").Value = 15874177
$ws.Range("K34").Value = 5403.0303
$ws.Range("L34").Value = 15874177
$ws.Range("M34").Value = -5201.0303
$ws.Range("N34").Value = -15874581

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H122").Value = 3378
$ws.Range("I122").Value = 2776.75
$ws.Range("J122").Value = 4340
$ws.Range("K122").Value = 8330.25
$ws.Range("L122").Value = 13020
$ws.Range("M122").Value = -5880.25
$ws.Range("N122").Value = -17920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 515.4
$ws.Range("I41").Value = 199
$ws.Range("K41").Value = 597
$ws.Range("M41").Value = -259

$ws.Range("H44").Value = 1127.0869
$ws.Range("I44").Value = 186.85715
$ws.Range("J44").Value = 1538.4375
$ws.Range("K44").Value = 560.5714499999999
$ws.Range("L44").Value = 4615.3125
$ws.Range("M44").Value = -162.5714499999999
$ws.Range("N44").Value = -5411.3125

$ws.Range("H62").Value = 5445.684
$ws.Range("I62").Value = 2253
$ws.Range("J62").Value = 6297.067
$ws.Range("K62").Value = 6759
$ws.Range("L62").Value = 18891.201
$ws.Range("M62").Value = -6073
$ws.Range("N62").Value = -20263.201

$ws.Range("H65").Value = 5445.684
$ws.Range("I65").Value = 2253
$ws.Range("J65").Value = 6297.067
$ws.Range("K65").Value = 20277
$ws.Range("L65").Value = 56673.603
$ws.Range("M65").Value = -16845
$ws.Range("N65").Value = -63537.603

$ws.Range("H68").Value = 1205.5918
$ws.Range("I68").Value = 943.45
$ws.Range("J68").Value = 1386.3793
$ws.Range("K68").Value = 2830.35
$ws.Range("L68").Value = 4159.1379
$ws.Range("M68").Value = -2019.35
$ws.Range("N68").Value = -5781.1379

$ws.Range("H71").Value = 1205.5918
$ws.Range("I71").Value = 943.45
$ws.Range("J71").Value = 1386.3793
$ws.Range("K71").Value = 8491.050000000001
$ws.Range("L71").Value = 12477.4137
$ws.Range("M71").Value = -4435.050000000001
$ws.Range("N71").Value = -20589.4137

$ws.Range("H131").Value = 739.4
$ws.Range("J131").Value = 817.0864
$ws.Range("L131").Value = 2451.2592
$ws.Range("N131").Value = -12531.2592

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 39800
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 39800
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 39800
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -40168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13000
$ws.Range("I40").Value = 18666.666
$ws.Range("J40").Value = 8750
$ws.Range("K40").Value = 18666.666
$ws.Range("L40").Value = 8750
$ws.Range("M40").Value = -18530.666
$ws.Range("N40").Value = -9022

$ws.Range("H94").Value = 49995
$ws.Range("J94").Value = 49995
$ws.Range("L94").Value = 49995
$ws.Range("N94").Value = -51347

$ws.Range("H100").Value = 1670
$ws.Range("I100").Value = 1603.75
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 1603.75
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -1062.75
$ws.Range("N100").Value = -3282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1500277.5
$ws.Range("J2").Value = 1500277.5
$ws.Range("L2").Value = 1500277.5
$ws.Range("N2").Value = -1500501.5

$ws.Range("H51").Value = 12500
$ws.Range("I51").Value = 12500
$ws.Range("K51").Value = 12500
$ws.Range("M51").Value = -11990

$ws.Range("H52").Value = 21000
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 22000
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 22000
$ws.Range("M52").Value = -19774
$ws.Range("N52").Value = -22452

$ws.Range("H100").Value = 1257.1428
$ws.Range("I100").Value = 750
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959

$ws.Range("H132").Value = 9844.071
$ws.Range("I132").Value = 15439.875
$ws.Range("J132").Value = 2383
$ws.Range("K132").Value = 46319.625
$ws.Range("L132").Value = 7149
$ws.Range("M132").Value = -43789.625
$ws.Range("N132").Value = -12209
